# Insert a new "quarter" column (2018-09-30) as column D on the FBSI
# worksheet, shifting the existing quarterly data (previously in D:K)
# one column to the right (now E:L), and fill in the new quarter's
# figures for the Income Statement, Balance Sheet and Cash Flow
# Statement blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column before D; Excel shifts D:K -> E:L automatically,
# carrying values/formats along with the cells being pushed.
$ws.Columns("D:D").Insert()

# The freshly inserted column D is blank and picks up a default format.
# Copy the per-row number format (date vs. numeric vs. blank) from the
# cells now sitting in column E (which used to be column D) back onto
# column D so every new-quarter cell matches its row's existing style.
# Done in separate contiguous blocks so we don't touch the header/blank
# separator rows (5, 6, 36, 37, 78, 79) that never had a D:K cell range.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Period Ending headers (new quarter: 2018-09-30)
$ws.Range("D7").Value = 43373
$ws.Range("D38").Value = 43373
$ws.Range("D80").Value = 43373

# Income Statement
$ws.Range("D8").Value = 11100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = "NA"
$ws.Range("D17").Value = 2200
$ws.Range("D18").Value = 8900
$ws.Range("D20").Value = -6000
$ws.Range("D21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 3000
$ws.Range("D24").Value = 700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 2200
$ws.Range("D27").Value = 2200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 6000
$ws.Range("D33").Value = 2200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 2200

# Balance Sheet
$ws.Range("D41").Value = 25600
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = "NA"
$ws.Range("D49").Value = 2400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 354000
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = "NA"
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = "NA"
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 32700
$ws.Range("D77").Value = 0

# Cash Flow Statement
$ws.Range("D81").Value = 2200
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0
